# extraction données et labels DF dans un dictionnaire
# (et retrait d'une ligne parasite dans le fichier excel 2023 et 2024)
#
# The sheet "2023-DRH-Annuel" has a blank/parasitic row (row 2, between the
# header row and the first data row "DRH-01") left over from data entry.
# Remove it entirely so every row below shifts up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2023-DRH-Annuel")

# Make sure this is the active/selected sheet (it was already the active tab).
$ws.Activate()

# Delete the whole parasitic empty row - everything below (old rows 3-7,
# now rows 2-6) shifts up, keeping its own formatting/formulas intact.
$ws.Rows(2).Delete()

# Leave the selection where the user ended up after the edit.
$ws.Range("D10").Select()
